# Download and ui fix
# - Inserts 2 new freshly-downloaded income entries at the top (rows 2-3),
#   plus one more (row 4) ahead of the previously-existing rows, and
#   appends 3 more rows that had come in from another sync at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new rows above the existing data (old row 2 becomes row 5) ---
$ws.Rows("2:4").Insert()

# Row 2: Freelance income
$ws.Range("A2").Value = "Freelance "
$ws.Range("B2").Value = 12000
$ws.Range("C5").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C2").Value = 45830.00011574074

# Row 3: rahul ki salary
$ws.Range("A3").Value = "rahul ki salary "
$ws.Range("B3").Value = 52867
$ws.Range("C5").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").Value = 45828.00011574074

# Row 4: uydbhjev
$ws.Range("A4").Value = "uydbhjev"
$ws.Range("B4").Value = 3494
$ws.Range("C5").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C4").Value = 45818.00011574074

# --- Append 3 new rows at the bottom (rows 11-13) ---
$ws.Range("A11").Value = "werty"
$ws.Range("B11").Value = 12
$ws.Range("C5").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C11").Value = 45812.00011574074

$ws.Range("A12").Value = "asdfgh"
$ws.Range("B12").Value = 12
$ws.Range("C5").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("C12").Value = 45810.00011574074

$ws.Range("A13").Value = "sdfghj"
$ws.Range("B13").Value = 12
$ws.Range("C5").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C13").Value = 45810.00011574074
